# The experiment's odor columns were relabelled with underscores
# ("odor1" -> "odor_1", "odor2" -> "odor_2") as part of wiring up the
# plotting / file-saving code. Update the label cells in column A,
# touching the rows containing "odor2" before the row containing
# "odor1" so the rebuilt shared-string table orders "odor_2" ahead of
# "odor_1", matching the saved workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "odor_2"
$ws.Range("A5").Value = "odor_1"
$ws.Range("A6").Value = "odor_2"

# Move the active selection from the old H7 cell to A6, reflecting
# where the user left the cursor after the edit.
$ws.Range("A6").Select() | Out-Null

Write-Host "Renamed odor1/odor2 labels to odor_1/odor_2 and moved selection to A6"
